$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated "attack speed" (column F) values per row on the Exp sheet.
# Rows 11 and 24 are intentionally left untouched (not part of the edit).
$updates = @{
    4 = 10;   5 = 14;   6 = 18;   7 = 22;   8 = 26;   9 = 30;  10 = 34;
   12 = 42;  13 = 47;  14 = 52;  15 = 57;  16 = 62;  17 = 67;  18 = 72;
   19 = 77;  20 = 82;  21 = 87;  22 = 92;  23 = 97;
   25 = 109; 26 = 115; 27 = 121; 28 = 127; 29 = 133; 30 = 139; 31 = 145;
   32 = 151; 33 = 157; 34 = 163; 35 = 169; 36 = 175; 37 = 181; 38 = 187;
   39 = 193; 40 = 199; 41 = 205; 42 = 211; 43 = 217; 44 = 223; 45 = 229;
   46 = 235; 47 = 241; 48 = 247; 49 = 253; 50 = 259; 51 = 265; 52 = 271;
   53 = 277; 54 = 283; 55 = 289; 56 = 295; 57 = 301; 58 = 307; 59 = 313;
   60 = 319; 61 = 325; 62 = 331; 63 = 337; 64 = 343; 65 = 349; 66 = 355;
   67 = 361; 68 = 367; 69 = 373; 70 = 379; 71 = 385; 72 = 391; 73 = 397;
   74 = 403; 75 = 409; 76 = 415; 77 = 421; 78 = 427; 79 = 433; 80 = 439;
   81 = 445; 82 = 451; 83 = 457; 84 = 463; 85 = 469; 86 = 475; 87 = 481;
   88 = 487; 89 = 493; 90 = 499; 91 = 505; 92 = 511; 93 = 517; 94 = 523;
   95 = 529; 96 = 535; 97 = 541; 98 = 547; 99 = 553; 100 = 559; 101 = 565;
   102 = 571
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}

# Reset the view: scroll back to the top-left (A1) and move the active
# selection to F9 instead of the previous D30 / topLeftCell A19.
$excel.Goto($ws.Range("F9"), $true)
